# Apply cryptos list update (Tue Jun 25 21:15:50 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.826.86"
$ws.Range("E2").Value = "  +3.84%  "
$ws.Range("D3").Value = "3.408.26"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("D5").Value = "'577.24"
$ws.Range("E5").Value = "  +2.81%  "
$ws.Range("D6").Value = "'137.72"
$ws.Range("E6").Value = "  +6.44%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.407.90"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "'7.50"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +8.44%  "
$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  +5.62%  "
$ws.Range("D13").Value = "3.995.44"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("E15").Value = "  +7.10%  "
$ws.Range("D16").Value = "3.407.72"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "'25.44"
$ws.Range("E17").Value = "  +4.40%  "
$ws.Range("D18").Value = "61.875.11"
$ws.Range("E18").Value = "  +3.53%  "
$ws.Range("D19").Value = "'14.15"
$ws.Range("E19").Value = "  +6.53%  "
$ws.Range("E20").Value = "  +4.69%  "
$ws.Range("D21").Value = "'9.48"
$ws.Range("E21").Value = "  +5.77%  "
$ws.Range("D22").Value = "'389.43"
$ws.Range("E22").Value = "  +10.82%  "
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").Value = "3.546.05"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "'0.0000128"
$ws.Range("E25").Value = "  +16.81%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'71.60"
$ws.Range("E27").Value = "  +4.39%  "
$ws.Range("D28").Value = "'7.70"
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").Value = "'1.58"
$ws.Range("E29").Value = "  +6.81%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  +5.97%  "
$ws.Range("E32").Value = "  +5.57%  "
$ws.Range("D33").Value = "'2.17"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("D34").Value = "3.439.15"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'23.57"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("D37").Value = "'5.45"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").Value = "'6.99"
$ws.Range("E38").Value = "  +3.16%  "
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("D40").Value = "'164.31"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "'0.0788"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("E42").Value = "  +14.05%  "
$ws.Range("D43").Value = "'0.789"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'41.70"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'24.96"
$ws.Range("E48").Value = "  +7.24%  "
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("D50").Value = "'23.04"
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("D51").Value = "2.378.75"
$ws.Range("E51").Value = "  +10.84%  "
